$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "From"
$ws.Range("C1").Value = "To"
$ws.Range("D1").Value = "Activity"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Status"
$ws.Range("A1").Value = "Date"

$ws.Columns.Item(1).ColumnWidth = 4.33
$ws.Columns.Item(5).ColumnWidth = 10.33

$ws.Range("B10").Select()
